$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like numbers (e.g. "316.22") need to be forced
# to remain plain text, matching the original inline-string cell type used for
# the "Price" column. We temporarily apply a text number format, assign the
# value, then restore the default "Normal" style so no stray formatting is left
# behind on the cell.
$textForcedAddresses = @(
    "D5",
    "D6",
    "D10",
    "D12",
    "D17",
    "D19",
    "D20",
    "D22",
    "D23",
    "D24",
    "D26",
    "D29",
    "D30",
    "D32",
    "D36",
    "D37",
    "D41",
    "D42",
    "D43",
    "D47",
    "D48",
    "D49",
    "D50"
)
foreach ($addr in $textForcedAddresses) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = "316.22"
$ws.Range("D6").Value = "96.34"
$ws.Range("D10").Value = "35.54"
$ws.Range("D12").Value = "7.49"
$ws.Range("D17").Value = "0.847"
$ws.Range("D19").Value = "6.81"
$ws.Range("D20").Value = "12.71"
$ws.Range("D22").Value = "69.55"
$ws.Range("D23").Value = "250.50"
$ws.Range("D24").Value = "2.94"
$ws.Range("D26").Value = "26.36"
$ws.Range("D29").Value = "40.44"
$ws.Range("D30").Value = "10.35"
$ws.Range("D32").Value = "157.59"
$ws.Range("D36").Value = "18.85"
$ws.Range("D37").Value = "0.0782"
$ws.Range("D41").Value = "22.34"
$ws.Range("D42").Value = "3.81"
$ws.Range("D43").Value = "0.0304"
$ws.Range("D47").Value = "9.03"
$ws.Range("D48").Value = "84.17"
$ws.Range("D49").Value = "105.94"
$ws.Range("D50").Value = "74.98"

foreach ($addr in $textForcedAddresses) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining cells are unambiguous text already (coin names, links, and the
# percentage strings), so a direct assignment is sufficient.
$ws.Range("D2").Value = "42.610.29"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "2.524.09"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("E13").Value = "  -2.93%  "
$ws.Range("D14").Value = "2.912.64"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").Value = "2.503.25"
$ws.Range("E15").Value = "  -2.05%  "
$ws.Range("E16").Value = "  -3.07%  "
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "42.709.39"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("E19").Value = "  +3.87%  "
$ws.Range("E20").Value = "  -4.08%  "
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("E22").Value = "  -2.52%  "
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -1.96%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("E29").Value = "  +3.69%  "
$ws.Range("E30").Value = "  +2.78%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("E33").Value = "  +2.47%  "
$ws.Range("E34").Value = "  +3.98%  "
$ws.Range("E35").Value = "  -1.04%  "
$ws.Range("E36").Value = "  -3.21%  "
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("E40").Value = "  +8.78%  "
$ws.Range("E41").Value = "  -6.45%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E42").Value = "  -1.61%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").Value = "2.024.53"
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("E47").Value = "  +2.32%  "
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("E49").Value = "  +4.05%  "
$ws.Range("E50").Value = "  +2.62%  "
$ws.Range("D51").Value = "2.767.68"
$ws.Range("E51").Value = "  +0.49%  "
